# Apply cryptos list update (Thu Jun  6 09:37:28 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '70.953.66'
$ws.Cells.Item(2, 5).Value = '  -0.04%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.850.66'
$ws.Cells.Item(3, 5).Value = '  +1.35%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.10%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '704.51'
$ws.Cells.Item(5, 5).Value = '  +1.24%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '173.06'
$ws.Cells.Item(6, 5).Value = '  -0.35%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '3.849.05'
$ws.Cells.Item(7, 5).Value = '  +1.38%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.01%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.83%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -0.54%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '7.37'
$ws.Cells.Item(11, 5).Value = '  -1.61%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -0.47%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -1.94%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '36.60'
$ws.Cells.Item(14, 5).Value = '  +0.68%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '4.498.94'
$ws.Cells.Item(15, 5).Value = '  +1.29%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.861.92'
$ws.Cells.Item(16, 5).Value = '  +1.59%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '70.993.99'
$ws.Cells.Item(17, 5).Value = '  +0.00%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  -0.17%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.69%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.39'
$ws.Cells.Item(20, 5).Value = '  -2.75%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -3.92%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '492.87'
$ws.Cells.Item(22, 5).Value = '  +1.69%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.717'
$ws.Cells.Item(23, 5).Value = '  +0.26%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '85.11'
$ws.Cells.Item(24, 5).Value = '  +1.01%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +1.01%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'RenderToken'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '10.58'
$ws.Cells.Item(26, 5).Value = '  +1.00%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '12.16'
$ws.Cells.Item(27, 5).Value = '  -2.25%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -2.08%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '3.18'
$ws.Cells.Item(29, 5).Value = '  +4.66%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.03%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '7.49'
$ws.Cells.Item(31, 5).Value = '  -0.36%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -0.84%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '29.47'
$ws.Cells.Item(33, 5).Value = '  -0.53%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.12%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.805.59'
$ws.Cells.Item(35, 5).Value = '  +1.46%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '9.17'
$ws.Cells.Item(36, 5).Value = '  -0.85%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.01%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +0.50%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +6.46%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.08'
$ws.Cells.Item(40, 5).Value = '  +1.65%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +6.44%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.33'
$ws.Cells.Item(42, 5).Value = '  -5.35%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +0.07%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '163.24'
$ws.Cells.Item(45, 5).Value = '  -0.22%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.000308'
$ws.Cells.Item(46, 5).Value = '  -5.68%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '48.74'
$ws.Cells.Item(47, 5).Value = '  -1.21%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'ONDO'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.39'
$ws.Cells.Item(48, 5).Value = '  +0.63%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Bittensor'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '415.40'
$ws.Cells.Item(49, 5).Value = '  +3.32%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'TheGraph'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.299'
$ws.Cells.Item(50, 5).Value = '  -0.76%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Cosmos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '8.62'
$ws.Cells.Item(51, 5).Value = '  +0.66%  '
